{"js": "// Update the ANOVA results table with the re-analyzed (subjective video\n// rating) statistics. The header row is left untouched; every data row's\n// Sum of Squares / Mean Square / F-Value / p-Value / Significance / Num df /\n// Den df columns are refreshed to the new numbers.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Leading spaces in some of the new values are significant (they are how the\n// source right-aligned the numbers as plain text), so keep them exactly.\nconst updates = {\n  1: [\"  7.22\", \" 7.22\", \"12.18\", \"0.001\", \"***\", \"1\", \"462\"],\n  2: [\"1018.70\", \"339.57\", \"572.69\", \"0.000\", \"***\", \"3\", \"462\"],\n  3: [\"  0.35\", \" 0.35\", \" 0.59\", \"0.443\", \"\", \"1\", \"462\"],\n  4: [\"  0.60\", \" 0.60\", \" 1.01\", \"0.315\", \"\", \"1\", \"462\"],\n  5: [\" 19.23\", \" 6.41\", \"10.81\", \"0.000\", \"***\", \"3\", \"462\"],\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = updates[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the ANOVA results table with the re-analyzed (subjective video\n# rating) statistics. The header row is left untouched; every data row's\n# Sum of Squares / Mean Square / F-Value / p-Value / Significance / Num df /\n# Den df columns are refreshed to the new numbers.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Leading spaces in some of the new values are significant (they are how the\n# source right-aligned the numbers as plain text), so keep them exactly.\n$updates = @{\n  2 = @(\"  7.22\", \" 7.22\", \"12.18\", \"0.001\", \"***\", \"1\", \"462\")\n  3 = @(\"1018.70\", \"339.57\", \"572.69\", \"0.000\", \"***\", \"3\", \"462\")\n  4 = @(\"  0.35\", \" 0.35\", \" 0.59\", \"0.443\", \"\", \"1\", \"462\")\n  5 = @(\"  0.60\", \" 0.60\", \" 1.01\", \"0.315\", \"\", \"1\", \"462\")\n  6 = @(\" 19.23\", \" 6.41\", \"10.81\", \"0.000\", \"***\", \"3\", \"462\")\n}\n\nforeach ($row in $updates.Keys) {\n  $values = $updates[$row]\n  for ($col = 1; $col -le $values.Length; $col++) {\n    $t.Cell($row, $col).Range.Text = $values[$col - 1]\n  }\n}\n"}
